# assets_liabilities.xlsx update
# - Summary sheet: borrower name/income/totals refreshed from the latest
#   pipeline run.
# - Assets sheet: vehicle re-valued (Mid-range -> Premium) and liquid
#   assets balance refreshed; total recomputed.
# - Liabilities sheet: auto loan balance refreshed, a new "Personal Loans"
#   line item inserted, credit-card balance refreshed, total recomputed.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------
# Summary sheet
# ---------------------------------------------------------------
$wsS = $wb.Worksheets.Item("Summary")
$wsS.Range("B3").Value = "Faisal Al Mansoori"
$wsS.Range("B4").Value = 2210.64
$wsS.Range("B6").Value = 188617
$wsS.Range("B7").Value = 177836
$wsS.Range("B8").Value = 10781
$wsS.Range("B9").Value = 1.06

# ---------------------------------------------------------------
# Assets sheet
# ---------------------------------------------------------------
$wsA = $wb.Worksheets.Item("Assets")
$wsA.Range("B2").Value = "Premium Car"
$wsA.Range("C2").Value = 186798
$wsA.Range("C3").Value = 1819
$wsA.Range("C4").Value = 188617

# ---------------------------------------------------------------
# Liabilities sheet
# ---------------------------------------------------------------
$wsL = $wb.Worksheets.Item("Liabilities")

# Refresh the Auto Loans row
$wsL.Range("C2").Value = 112079
$wsL.Range("D2").Value = 2335
$wsL.Range("E2").Value = 4

# Insert a new "Personal Loans" row above the Credit Cards row, copying
# the formatting of the row above it (Auto Loans) so the banded style
# carries through.
$wsL.Rows.Item(3).Insert()
$wsL.Range("A2:E2").Copy()
$wsL.Range("A3:E3").PasteSpecial(-4122)

$wsL.Range("A3").Value = "Personal Loans"
$wsL.Range("B3").Value = "Personal Loan"
$wsL.Range("C3").Value = 29784
$wsL.Range("D3").Value = 620
$wsL.Range("E3").Value = 4

# Credit Cards row (shifted down to row 4 by the insert above)
$wsL.Range("C4").Value = 35973
$wsL.Range("D4").Value = 1799

# TOTAL LIABILITIES row (shifted down to row 5)
$wsL.Range("C5").Value = 177836
